$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their text type (avoid Excel auto-numeric coercion)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.732.10'
$ws.Range('E2').Value = '  -2.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.557.73'
$ws.Range('E3').Value = '  -3.39%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '514.59'
$ws.Range('E5').Value = '  -2.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.57'
$ws.Range('E6').Value = '  -5.63%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.558'
$ws.Range('E8').Value = '  -2.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.556.88'
$ws.Range('E9').Value = '  -3.87%  '
$ws.Range('E10').Value = '  -3.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0982'
$ws.Range('E11').Value = '  -5.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.323'
$ws.Range('E12').Value = '  -4.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.011.91'
$ws.Range('E14').Value = '  -3.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.725.45'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.92'
$ws.Range('E16').Value = '  -4.42%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.570.87'
$ws.Range('E17').Value = '  -2.43%  '
$ws.Range('E18').Value = '  -5.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '331.10'
$ws.Range('E19').Value = '  -3.44%  '
$ws.Range('E20').Value = '  -4.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.99'
$ws.Range('E21').Value = '  -6.08%  '
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.37'
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -5.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.87'
$ws.Range('E28').Value = '  -4.81%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0699'
$ws.Range('E30').Value = '  -12.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.88'
$ws.Range('E31').Value = '  -8.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.77'
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('E33').Value = '  -4.17%  '
$ws.Range('E34').Value = '  -2.96%  '
$ws.Range('E35').Value = '  -7.44%  '
$ws.Range('E36').Value = '  -7.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '35.86'
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.822'
$ws.Range('E38').Value = '  -4.64%  '
$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.815'
$ws.Range('E39').Value = '  -5.85%  '
$ws.Range('E40').Value = '  -5.96%  '
$ws.Range('E41').Value = '  -5.08%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '268.51'
$ws.Range('E43').Value = '  -0.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.70'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0935'
$ws.Range('E45').Value = '  -4.33%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.578'
$ws.Range('E46').Value = '  -3.93%  '
$ws.Range('E47').Value = '  -4.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.952.97'
$ws.Range('E48').Value = '  -4.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.17'
$ws.Range('E49').Value = '  -6.31%  '
$ws.Range('E50').Value = '  -5.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.35'
$ws.Range('E51').Value = '  -8.86%  '
